$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "71.562.66"
$ws.Range("D2").Style = $style
$ws.Range("E2").Value = "  -1.99%  "

$style = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.883.74"
$ws.Range("D3").Style = $style
$ws.Range("E3").Value = "  -2.93%  "

$style = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = $style
$ws.Range("E4").Value = "  +0.01%  "

$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.29"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  -2.28%  "

$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.15"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  +6.23%  "

$style = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.670"
$ws.Range("D7").Style = $style
$ws.Range("E7").Value = "  -2.68%  "

$style = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = $style
$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("E9").Value = "  -1.51%  "

$style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.178"
$ws.Range("D10").Style = $style
$ws.Range("E10").Value = "  +4.84%  "

$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.13"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = "  -0.57%  "

$style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000324"
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = "  +0.86%  "

$style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.46"
$ws.Range("D13").Style = $style
$ws.Range("E13").Value = "  +2.72%  "

$style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.501.52"
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = "  -3.00%  "

$style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.07"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = "  +1.64%  "

$style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.873.92"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = "  -3.20%  "

$ws.Range("E17").Value = "  -1.85%  "

$ws.Range("E18").Value = "  -4.14%  "

$ws.Range("E19").Value = "  -2.27%  "

$style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.434.91"
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = "  -1.76%  "

$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "440.53"
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = "  -0.28%  "

$ws.Range("E22").Value = "  -2.48%  "

$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "94.39"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = "  -2.67%  "

$ws.Range("E24").Value = "  -4.36%  "

$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.91"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  -4.20%  "

$ws.Range("E26").Value = "  +3.03%  "

$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.06"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  -5.52%  "

$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.98"
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = "  +0.51%  "

$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.51"
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = "  -1.16%  "

$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.70"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = "  +11.87%  "

$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.25"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = "  -3.70%  "

$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.62"
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = "  -2.97%  "

$ws.Range("E33").Value = "  -3.97%  "

$ws.Range("B34").Value = "PEPE"
$ws.Range("C34").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0000102"
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = "  +10.75%  "

$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "47.75"
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = "  -1.19%  "

$style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "69.90"
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = "  -3.86%  "

$style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "634.67"
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = "  -2.75%  "

$style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.437"
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = "  -0.91%  "

$ws.Range("E39").Value = "  -0.47%  "

$ws.Range("E40").Value = "  +0.16%  "

$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.32"
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = "  -2.69%  "

$ws.Range("E42").Value = "  -0.19%  "

$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.90"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = "  +9.02%  "

$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.20"
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = "  +18.88%  "

$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0473"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = "  -3.91%  "

$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.26"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = "  -3.94%  "

$ws.Range("E47").Value = "  -3.89%  "

$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.91"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  -13.18%  "

$style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.914.65"
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = "  -0.43%  "

$style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000281"
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = "  +3.49%  "

$ws.Range("E51").Value = "  -5.29%  "
